# MAS Report 1A template formatting update
# - Apply currency number format to the "value" cells of the two summary tables
# - Clear the placeholder "0" values that were pre-filled in the template's input cells
# - Reset the active sheet's view (no frozen/scrolled topLeftCell, selection moved to C25)
# - Resize / reposition the workbook window (best effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFormat = """$""#,##0.00"

# 1. Apply the new currency ("$"#,##0.00) number format to the monetary "value" input cells.
foreach ($addr in @("B6", "B7", "B10", "B11")) {
    $ws.Range($addr).NumberFormat = $currencyFormat
}

# 2. Clear the placeholder 0 values left over in the template's input cells.
#    Cells B14:C14, B15:C15 ... B19:C19 and B22:C22 are merged, so the whole
#    merged range must be targeted for the clear to stick.
foreach ($addr in @("B6", "C6", "B7", "C7", "B10", "C10", "B11", "C11", "B14:C14", "B15:C15", "B16:C16", "B17:C17", "B18:C18", "B19:C19", "B22:C22")) {
    $ws.Range($addr).ClearContents()
}

# 3. Update the sheet view: drop the frozen/scrolled topLeftCell and move the selection to C25.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("C25").Select()

# 4. Best-effort update of the workbook window position/size.
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840

$wb.Save()
